$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEB-22")

# Row 36 - continuation note (set first so it claims the next shared-string slot)
$ws.Range("D36").Value = "whereas no merging is taking place during update, but, still error is throwing"

# Row 38 - continuation note
$ws.Range("D38").Value = "Hence, i have uninstalled lower version of openpyxl and it is work in progress."

# Row 35 - new log entry (row 23 in the daily log)
$ws.Range("A35").Value = 23

# Use an existing date cell as a template so the date keeps the workbook's
# established short-date style instead of minting a brand-new number format.
[void]$ws.Range("B33").Copy($ws.Range("B35"))
$ws.Range("B35").Value = 44621

$ws.Range("C35").Value = "RPA GSS"
$ws.Range("D35").Value = "1. In invoice generation task, while updating the master file from the csv file is throwing error :"
$ws.Range("E35").Value = 0.85
$ws.Range("F35").Value = "WIP"

# Row 37 - continuation note
$ws.Range("D37").Value = "2. After using lower version of openpyxl, the master file is getting updated  success, but somewhere a few formulas are getting  error or issue."

# Update the active view/selection to reflect the scrolled-down state
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
[void]$ws.Range("D40").Select()
